$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 6).Value = 1
$ws.Cells.Item(2, 7).Value = 14.75986366666667
$ws.Cells.Item(2, 8).Value = 44.279591
$ws.Cells.Item(2, 9).Value = 0.2069066005118911
$ws.Cells.Item(2, 10).Value = 0.2199725985531551
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 12).Value = 1
$ws.Cells.Item(2, 13).Value = 0.7727323333333332
$ws.Cells.Item(2, 14).Value = 2.318197
$ws.Cells.Item(2, 15).Value = 0.08417634149364144
$ws.Cells.Item(2, 16).Value = 0.09406687683971243
$ws.Cells.Item(2, 17).Value = 11.40542389082522
$ws.Cells.Item(2, 18).Value = 102.648815017427
$ws.Cells.Item(2, 19).Value = 0.01741664066197739
$ws.Cells.Item(2, 20).Value = 0.02069213533621114

# Row 3
$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 6).Value = 1
$ws.Cells.Item(3, 7).Value = 14.75986366666667
$ws.Cells.Item(3, 8).Value = 44.279591
$ws.Cells.Item(3, 9).Value = 0.2069066005118911
$ws.Cells.Item(3, 10).Value = 0.2199725985531551
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(3, 13).Value = 5.511558666666666
$ws.Cells.Item(3, 14).Value = 16.534676
$ws.Cells.Item(3, 15).Value = 0.6003926902945337
$ws.Cells.Item(3, 16).Value = 0.6709375134540114
$ws.Cells.Item(3, 17).Value = 81.3498545108351
$ws.Cells.Item(3, 18).Value = 732.1486905975158
$ws.Cells.Item(3, 19).Value = 0.1242252105210306
$ws.Cells.Item(3, 20).Value = 0.1475878683012713

# Row 4
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 6).Value = 1
$ws.Cells.Item(4, 7).Value = 14.75986366666667
$ws.Cells.Item(4, 8).Value = 44.279591
$ws.Cells.Item(4, 9).Value = 0.2069066005118911
$ws.Cells.Item(4, 10).Value = 0.2199725985531551
$ws.Cells.Item(4, 11).Value = 2
$ws.Cells.Item(4, 12).Value = 1
$ws.Cells.Item(4, 13).Value = 2.895632
$ws.Cells.Item(4, 14).Value = 5.791264
$ws.Cells.Item(4, 15).Value = 0.3154309682118249
$ws.Cells.Item(4, 16).Value = 0.2349956097062763
$ws.Cells.Item(4, 17).Value = 42.73913354883733
$ws.Cells.Item(4, 18).Value = 256.434801293024
$ws.Cells.Item(4, 19).Value = 0.06526474932888307
$ws.Cells.Item(4, 20).Value = 0.05169259491567262

# Row 5
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 10.98286
$ws.Cells.Item(5, 8).Value = 32.94858
$ws.Cells.Item(5, 9).Value = 0.1539598385065048
$ws.Cells.Item(5, 10).Value = 0.1636822878792289
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 12).Value = 1
$ws.Cells.Item(5, 13).Value = 0.7727323333333332
$ws.Cells.Item(5, 14).Value = 2.318197
$ws.Cells.Item(5, 15).Value = 0.08417634149364144
$ws.Cells.Item(5, 16).Value = 0.09406687683971243
$ws.Cells.Item(5, 17).Value = 8.486811034473332
$ws.Cells.Item(5, 18).Value = 76.38129931025999
$ws.Cells.Item(5, 19).Value = 0.01295977594242943
$ws.Cells.Item(5, 20).Value = 0.01539708161477778

# Row 6
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 6).Value = 1
$ws.Cells.Item(6, 7).Value = 10.98286
$ws.Cells.Item(6, 8).Value = 32.94858
$ws.Cells.Item(6, 9).Value = 0.1539598385065048
$ws.Cells.Item(6, 10).Value = 0.1636822878792289
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 12).Value = 1
$ws.Cells.Item(6, 13).Value = 5.511558666666666
$ws.Cells.Item(6, 14).Value = 16.534676
$ws.Cells.Item(6, 15).Value = 0.6003926902945337
$ws.Cells.Item(6, 16).Value = 0.6709375134540114
$ws.Cells.Item(6, 17).Value = 60.53267721778666
$ws.Cells.Item(6, 18).Value = 544.7940949600799
$ws.Cells.Item(6, 19).Value = 0.09243636163823235
$ws.Cells.Item(6, 20).Value = 0.1098205872261535

# Row 7
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 6).Value = 1
$ws.Cells.Item(7, 7).Value = 10.98286
$ws.Cells.Item(7, 8).Value = 32.94858
$ws.Cells.Item(7, 9).Value = 0.1539598385065048
$ws.Cells.Item(7, 10).Value = 0.1636822878792289
$ws.Cells.Item(7, 11).Value = 2
$ws.Cells.Item(7, 12).Value = 1
$ws.Cells.Item(7, 13).Value = 2.895632
$ws.Cells.Item(7, 14).Value = 5.791264
$ws.Cells.Item(7, 15).Value = 0.3154309682118249
$ws.Cells.Item(7, 16).Value = 0.2349956097062763
$ws.Cells.Item(7, 17).Value = 31.80232086752
$ws.Cells.Item(7, 18).Value = 190.81392520512
$ws.Cells.Item(7, 19).Value = 0.04856370092584301
$ws.Cells.Item(7, 20).Value = 0.03846461903829763

# Row 8
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 6).Value = 1
$ws.Cells.Item(8, 7).Value = 15.43916166666667
$ws.Cells.Item(8, 8).Value = 46.317485
$ws.Cells.Item(8, 9).Value = 0.2164291301970361
$ws.Cells.Item(8, 10).Value = 0.2300964689104012
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(8, 12).Value = 1
$ws.Cells.Item(8, 13).Value = 0.7727323333333332
$ws.Cells.Item(8, 14).Value = 2.318197
$ws.Cells.Item(8, 15).Value = 0.08417634149364144
$ws.Cells.Item(8, 16).Value = 0.09406687683971243
$ws.Cells.Item(8, 17).Value = 11.93033941939389
$ws.Cells.Item(8, 18).Value = 107.373054774545
$ws.Cells.Item(8, 19).Value = 0.01821821237263749
$ws.Cells.Item(8, 20).Value = 0.02164445620224743

# Row 9
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 6).Value = 1
$ws.Cells.Item(9, 7).Value = 15.43916166666667
$ws.Cells.Item(9, 8).Value = 46.317485
$ws.Cells.Item(9, 9).Value = 0.2164291301970361
$ws.Cells.Item(9, 10).Value = 0.2300964689104012
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 12).Value = 1
$ws.Cells.Item(9, 13).Value = 5.511558666666666
$ws.Cells.Item(9, 14).Value = 16.534676
$ws.Cells.Item(9, 15).Value = 0.6003926902945337
$ws.Cells.Item(9, 16).Value = 0.6709375134540114
$ws.Cells.Item(9, 17).Value = 85.09384528998443
$ws.Cells.Item(9, 18).Value = 765.8446076098598
$ws.Cells.Item(9, 19).Value = 0.1299424677371044
$ws.Cells.Item(9, 20).Value = 0.1543803527052929

# Row 10
$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 6).Value = 1
$ws.Cells.Item(10, 7).Value = 15.43916166666667
$ws.Cells.Item(10, 8).Value = 46.317485
$ws.Cells.Item(10, 9).Value = 0.2164291301970361
$ws.Cells.Item(10, 10).Value = 0.2300964689104012
$ws.Cells.Item(10, 11).Value = 2
$ws.Cells.Item(10, 12).Value = 1
$ws.Cells.Item(10, 13).Value = 2.895632
$ws.Cells.Item(10, 14).Value = 5.791264
$ws.Cells.Item(10, 15).Value = 0.3154309682118249
$ws.Cells.Item(10, 16).Value = 0.2349956097062763
$ws.Cells.Item(10, 17).Value = 44.70613057517333
$ws.Cells.Item(10, 18).Value = 268.23678345104
$ws.Cells.Item(10, 19).Value = 0.06826845008729419
$ws.Cells.Item(10, 20).Value = 0.05407166000286098

# Row 11
$ws.Cells.Item(11, 5).Value = 3
$ws.Cells.Item(11, 6).Value = 1
$ws.Cells.Item(11, 7).Value = 17.442302
$ws.Cells.Item(11, 8).Value = 52.326906
$ws.Cells.Item(11, 9).Value = 0.2445095356857582
$ws.Cells.Item(11, 10).Value = 0.2599501311352827
$ws.Cells.Item(11, 11).Value = 3
$ws.Cells.Item(11, 12).Value = 1
$ws.Cells.Item(11, 13).Value = 0.7727323333333332
$ws.Cells.Item(11, 14).Value = 2.318197
$ws.Cells.Item(11, 15).Value = 0.08417634149364144
$ws.Cells.Item(11, 16).Value = 0.09406687683971243
$ws.Cells.Item(11, 17).Value = 13.47823072316467
$ws.Cells.Item(11, 18).Value = 121.304076508482
$ws.Cells.Item(11, 19).Value = 0.02058191817433609
$ws.Cells.Item(11, 20).Value = 0.02445269696996973

# Row 12
$ws.Cells.Item(12, 5).Value = 3
$ws.Cells.Item(12, 6).Value = 1
$ws.Cells.Item(12, 7).Value = 17.442302
$ws.Cells.Item(12, 8).Value = 52.326906
$ws.Cells.Item(12, 9).Value = 0.2445095356857582
$ws.Cells.Item(12, 10).Value = 0.2599501311352827
$ws.Cells.Item(12, 11).Value = 3
$ws.Cells.Item(12, 12).Value = 1
$ws.Cells.Item(12, 13).Value = 5.511558666666666
$ws.Cells.Item(12, 14).Value = 16.534676
$ws.Cells.Item(12, 15).Value = 0.6003926902945337
$ws.Cells.Item(12, 16).Value = 0.6709375134540114
$ws.Cells.Item(12, 17).Value = 96.13427075471732
$ws.Cells.Item(12, 18).Value = 865.2084367924559
$ws.Cells.Item(12, 19).Value = 0.1468017379330396
$ws.Cells.Item(12, 20).Value = 0.1744102946059508

# Row 13
$ws.Cells.Item(13, 5).Value = 3
$ws.Cells.Item(13, 6).Value = 1
$ws.Cells.Item(13, 7).Value = 17.442302
$ws.Cells.Item(13, 8).Value = 52.326906
$ws.Cells.Item(13, 9).Value = 0.2445095356857582
$ws.Cells.Item(13, 10).Value = 0.2599501311352827
$ws.Cells.Item(13, 11).Value = 2
$ws.Cells.Item(13, 12).Value = 1
$ws.Cells.Item(13, 13).Value = 2.895632
$ws.Cells.Item(13, 14).Value = 5.791264
$ws.Cells.Item(13, 15).Value = 0.3154309682118249
$ws.Cells.Item(13, 16).Value = 0.2349956097062763
$ws.Cells.Item(13, 17).Value = 50.506487824864
$ws.Cells.Item(13, 18).Value = 303.038926949184
$ws.Cells.Item(13, 19).Value = 0.07712587957838246
$ws.Cells.Item(13, 20).Value = 0.06108713955936222

# Row 14
$ws.Cells.Item(14, 5).Value = 2
$ws.Cells.Item(14, 6).Value = 1
$ws.Cells.Item(14, 7).Value = 12.711689
$ws.Cells.Item(14, 8).Value = 25.423378
$ws.Cells.Item(14, 9).Value = 0.1781948950988097
$ws.Cells.Item(14, 10).Value = 0.1262985135219319
$ws.Cells.Item(14, 11).Value = 3
$ws.Cells.Item(14, 12).Value = 1
$ws.Cells.Item(14, 13).Value = 0.7727323333333332
$ws.Cells.Item(14, 14).Value = 2.318197
$ws.Cells.Item(14, 15).Value = 0.08417634149364144
$ws.Cells.Item(14, 16).Value = 0.09406687683971243
$ws.Cells.Item(14, 17).Value = 9.822733101577665
$ws.Cells.Item(14, 18).Value = 58.93639860946599
$ws.Cells.Item(14, 19).Value = 0.01499979434226102
$ws.Cells.Item(14, 20).Value = 0.01188050671650633

# Row 15
$ws.Cells.Item(15, 5).Value = 2
$ws.Cells.Item(15, 6).Value = 1
$ws.Cells.Item(15, 7).Value = 12.711689
$ws.Cells.Item(15, 8).Value = 25.423378
$ws.Cells.Item(15, 9).Value = 0.1781948950988097
$ws.Cells.Item(15, 10).Value = 0.1262985135219319
$ws.Cells.Item(15, 11).Value = 3
$ws.Cells.Item(15, 12).Value = 1
$ws.Cells.Item(15, 13).Value = 5.511558666666666
$ws.Cells.Item(15, 14).Value = 16.534676
$ws.Cells.Item(15, 15).Value = 0.6003926902945337
$ws.Cells.Item(15, 16).Value = 0.6709375134540114
$ws.Cells.Item(15, 17).Value = 70.06121967592132
$ws.Cells.Item(15, 18).Value = 420.3673180555279
$ws.Cells.Item(15, 19).Value = 0.1069869124651266
$ws.Cells.Item(15, 20).Value = 0.08473841061534283

# Row 16
$ws.Cells.Item(16, 5).Value = 2
$ws.Cells.Item(16, 6).Value = 1
$ws.Cells.Item(16, 7).Value = 12.711689
$ws.Cells.Item(16, 8).Value = 25.423378
$ws.Cells.Item(16, 9).Value = 0.1781948950988097
$ws.Cells.Item(16, 10).Value = 0.1262985135219319
$ws.Cells.Item(16, 11).Value = 2
$ws.Cells.Item(16, 12).Value = 1
$ws.Cells.Item(16, 13).Value = 2.895632
$ws.Cells.Item(16, 14).Value = 5.791264
$ws.Cells.Item(16, 15).Value = 0.3154309682118249
$ws.Cells.Item(16, 16).Value = 0.2349956097062763
$ws.Cells.Item(16, 17).Value = 44.70613057517333
$ws.Cells.Item(16, 18).Value = 147.233493769792
$ws.Cells.Item(16, 19).Value = 0.05620818829142213
$ws.Cells.Item(16, 20).Value = 0.02967959619008277

